$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V3.2-evalboards")

# Row 13 gains a new BOM entry (button / FSM4JH). Replicate the look & feel
# of the existing rows by copying the number/value-column formatting from
# row 12 before writing the new values.
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new values (write A13/E13 before C13/D13 so the shared
# string table order matches: button, link, FSM4JH, 506-FSM4JH).
$ws.Range("A13").Value = "button"
$ws.Range("E13").Value = "https://www.mouser.de/ProductDetail/TE-Connectivity-Alcoswitch/FSM4JH?qs=g%252BEszo6zu8OwVWrHD2r3Rw=="
$ws.Range("C13").Value = "FSM4JH"
$ws.Range("D13").Value = "506-FSM4JH"

# Match the saved selection state (active cell moved to D13).
$ws.Range("D13").Select()
